$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Code sets")

# Remember column B's width before insertion so the new column can match it
$colBWidth = $ws.Columns.Item(2).ColumnWidth

# Insert a new column before the existing "Occurrences required" column (C)
$ws.Columns.Item(3).Insert()

# Give the newly inserted column the same width as column B (matches Excel's
# typical behavior of carrying over neighboring column formatting)
$ws.Columns.Item(3).ColumnWidth = $colBWidth

# Header for the new column
$ws.Range("C1").Value = "Incident required"

# Rows that do NOT require an incident (2 occurrences needed instead)
$nRows = @(5,6,7,10,11,12)

for ($r = 2; $r -le 28; $r++) {
    if ($nRows -contains $r) {
        $ws.Range("C$r").Value = "N"
    } else {
        $ws.Range("C$r").Value = "Y"
    }
}

# Fill in the previously-blank "Occurrences required" values for rows 13 & 14
$ws.Range("D13").Value = 1
$ws.Range("D14").Value = 1

# Correct "Occurrences required" value for row 3 (was 1, now 2)
$ws.Range("D3").Value = 2
